$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the header row with the four values (matches commit's new sheetData)
$ws.Range("A1").Value = "wee"
$ws.Range("B1").Value = "are "
$ws.Range("C1").Value = "THE "
$ws.Range("D1").Value = "BEST"

# After typing the last entry, Excel leaves the selection on the next row;
# reflect that so the saved view matches (selection activeCell="A2").
[void]$ws.Range("A2").Select()
